$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cell D1 "stfips", matching style of existing header cells
$ws.Range("D1").Value = "stfips"
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122) # xlPasteFormats

# Update data values (A:C) with the refreshed dataset and add stfips "05" in column D
$data = @(
    @(742592, 380871, 361721),
    @(277047, 140733, 136314),
    @(393214, 196609, 196605),
    @(376107, 188040, 188067),
    @(359672, 178734, 180938),
    @(385089, 187380, 197709),
    @(300153, 140228, 159925),
    @(192017, 80583, 111434),
    @(3025891, 1493178, 1532713)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
    $d = $ws.Cells.Item($row, 4)
    $d.NumberFormat = "@"
    $d.Value = "05"
    $d.Style = "Normal"
}

$wb.Save()
